$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.39510000000003
$ws.Range("C5").Value = -13.93169999999999
$ws.Range("D7").Value = -7.394199999999995
$ws.Range("C9").Value = -12.10100000000001
$ws.Range("C11").Value = -12.8593
$ws.Range("D11").Value = -8.192000000000002
$ws.Range("A21").Value = -21.13500000000001
$ws.Range("C21").Value = -10.57279999999999
$ws.Range("D21").Value = -7.171199999999998
$ws.Range("A23").Value = -21.41550000000002
$ws.Range("A25").Value = -22.56270000000003
